# Generate Report for Handback
# Update status text and append error detail for the handback transform
# failure on the 48b1637f-...-947 row in both locale sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"
$zhCnError  = "Handback file name: 1pfdtgmw.wwa is different with handoff file name: 48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.zh-cn."
$deDeError  = "Handback file name: 1pfdtgmw.wwa is different with handoff file name: 48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.de-de."

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the status for the 48b1637f row (row 3) in all three sheets.
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText
$zhcn.Range("C3").Value = $statusText
$dede.Range("C3").Value = $statusText

# Add the error detail message to column K row 3 of each locale sheet.
$zhcn.Range("K3").Value = $zhCnError
$dede.Range("K3").Value = $deDeError
